# "Generate Report for Handback" — refresh the localization-status report
# after a handback: update status text, record the handback file/datetime
# for each locale, and widen a few columns so the new (longer) text fits.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$mdFile  = "c9a46d89-45ba-4db9-838f-7659a4255cb9.md"
$mdUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/388d589be37742e13bca914744a031051fee7a61/e2e/c9a46d89-45ba-4db9-838f-7659a4255cb9.md"
$status  = "Handed back: in sync with en-US"

# --- Overview sheet: status shown per-locale -------------------------------
$wsOverview.Range("E2").Value = $status
$wsOverview.Range("F2").Value = $status

# --- zh-cn sheet -------------------------------------------------------------
$wsZhCn.Range("C2").Value = $status
$wsZhCn.Range("I2").Value = $mdFile
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFile)
$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = 15570276
$wsZhCn.Range("J2").Value = "c9a46d89-45ba-4db9-838f-7659a4255cb9.94665bda437ee677dc4f3e8b9f53d435807ebe71.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-19 02:55:44"

# --- de-de sheet -------------------------------------------------------------
$wsDeDe.Range("C2").Value = $status
$wsDeDe.Range("I2").Value = $mdFile
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFile)
$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = 15570276
$wsDeDe.Range("J2").Value = "c9a46d89-45ba-4db9-838f-7659a4255cb9.94665bda437ee677dc4f3e8b9f53d435807ebe71.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-19 02:55:51"

# --- Column widths: widen columns now holding the longer status text and
# the newly-populated handoff/handback filename columns.
$wsOverview.Columns.Item(5).ColumnWidth = 29.144371396019366
$wsOverview.Columns.Item(6).ColumnWidth = 29.144371396019366

$wsZhCn.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDeDe.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
